# Applies the "Lector ocular terminado, aplicacion terminada" commit:
#  - Renames every "TAB." worksheet to its uppercase "TAB. ..." caption.
#  - Moves the active/selected tab from "casa" (TAB. CASA) to
#    "conceptos" (TAB. CONCEPTOS).

$wb = $excel.ActiveWorkbook

$names = @(
    "TAB. INICIAL",
    "TAB. RÁPIDO",
    "TAB. VERBOS",
    "TAB. VERBOS2",
    "TAB. COMIDA",
    "TAB. COMIDA2",
    "TAB. OBXECTOS",
    "TAB. PERSOAS",
    "TAB. LUGARES",
    "TAB. TRANSPORTE",
    "TAB. CASA",
    "TAB. ANIMAIS",
    "TAB. CORPO",
    "TAB. CONCEPTOS"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $names[$i - 1]
}

# The last worksheet ("TAB. CONCEPTOS") becomes the active / selected tab,
# replacing the previously active "TAB. CASA" sheet.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
